# Fix typos that were pointed out in comments received, and leave the
# "BothFilter" sheet as the active tab/selection, as was the case when the
# workbook was saved after making these edits.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("BothFilter")

# Fix "targetting" -> "targeting"
$ws.Range("A4").Value = "Remove trips targeting offshore species"

# Fix "co-occurand" -> "co-occur and" (both occurrences of this filter step)
$ws.Range("A7").Value = "Remove species that never co-occur and  not present in at least 1% of all"
$ws.Range("A13").Value = "Remove species that never co-occur and  not present in at least 1% of all"

# Make BothFilter the active sheet, with A13 selected
$ws.Activate()
$ws.Range("A13").Select()
